$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "61.947.48"
$ws.Range("E2").Value = "  +1.56%  "
Set-TextValue "D3" "3.418.91"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue "D5" "578.77"
$ws.Range("E5").Value = "  +1.32%  "
Set-TextValue "D6" "144.18"
$ws.Range("E6").Value = "  +2.39%  "
Set-TextValue "D8" "0.474"
$ws.Range("E8").Value = "  +0.08%  "
$ws.Range("E9").Value = "  -0.71%  "
Set-TextValue "D10" "0.123"
$ws.Range("E10").Value = "  +0.94%  "
Set-TextValue "D11" "0.386"
$ws.Range("E11").Value = "  -0.29%  "
Set-TextValue "D12" "4.006.14"
$ws.Range("E12").Value = "  +1.15%  "
Set-TextValue "D13" "28.52"
$ws.Range("E13").Value = "  +2.40%  "
$ws.Range("E14").Value = "  -0.62%  "
Set-TextValue "D15" "3.418.68"
$ws.Range("E15").Value = "  +1.20%  "
Set-TextValue "D16" "0.0000170"
$ws.Range("E16").Value = "  +0.10%  "
Set-TextValue "D17" "62.042.60"
$ws.Range("E17").Value = "  +1.54%  "
Set-TextValue "D18" "6.18"
$ws.Range("E18").Value = "  +1.12%  "
Set-TextValue "D19" "14.04"
$ws.Range("E19").Value = "  +3.21%  "
Set-TextValue "D20" "9.18"
$ws.Range("E20").Value = "  +3.27%  "
Set-TextValue "D21" "392.19"
$ws.Range("E21").Value = "  +2.33%  "
Set-TextValue "D22" "74.79"
$ws.Range("E22").Value = "  -2.04%  "
Set-TextValue "D23" "0.554"
$ws.Range("E23").Value = "  +0.44%  "
$ws.Range("E24").Value = "  +0.13%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
Set-TextValue "D25" "3.564.30"
$ws.Range("E25").Value = "  +1.24%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D26" "0.0000115"
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  +1.08%  "
Set-TextValue "D28" "7.50"
$ws.Range("E28").Value = "  +3.95%  "
$ws.Range("E29").Value = "  -0.11%  "
Set-TextValue "D30" "8.01"
$ws.Range("E30").Value = "  +0.77%  "
$ws.Range("E31").Value = "  +0.57%  "
Set-TextValue "D32" "1.40"
$ws.Range("E32").Value = "  +2.58%  "
$ws.Range("E33").Value = "  +0.00%  "
Set-TextValue "D34" "23.59"
Set-TextValue "D35" "5.29"
$ws.Range("E35").Value = "  +6.20%  "
Set-TextValue "D36" "6.97"
$ws.Range("E36").Value = "  +0.57%  "
Set-TextValue "D37" "167.78"
$ws.Range("E37").Value = "  +1.31%  "
Set-TextValue "D38" "1.53"
$ws.Range("E38").Value = "  +5.26%  "
Set-TextValue "D39" "3.450.57"
$ws.Range("E39").Value = "  +0.97%  "
Set-TextValue "D40" "29.21"
$ws.Range("E40").Value = "  +10.64%  "
Set-TextValue "D41" "0.0754"
$ws.Range("E41").Value = "  -1.62%  "
Set-TextValue "D42" "0.787"
$ws.Range("E42").Value = "  +1.23%  "
Set-TextValue "D43" "4.43"
$ws.Range("E43").Value = "  +1.81%  "
Set-TextValue "D44" "1.67"
$ws.Range("E44").Value = "  +1.47%  "
$ws.Range("E45").Value = "  +4.13%  "
Set-TextValue "D46" "2.513.19"
Set-TextValue "D47" "22.92"
$ws.Range("E47").Value = "  +0.36%  "
Set-TextValue "D48" "6.66"
$ws.Range("E48").Value = "  +0.36%  "
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("E50").Value = "  +1.16%  "
Set-TextValue "D51" "2.11"
$ws.Range("E51").Value = "  -0.62%  "
